# 📊 Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Drivers" table: clear the stale driver-vintage date for the
# 22.150.3.1 adapter driver row (E12) - no longer reported this week.
$ws.Range("E12").Value = ""

# Update total-samples count for the 22.100.1.1 adapter driver row (B14).
$ws.Range("B14").Value = 265400
